# Rename "Sheet1" to "template" and keep the Print_Area defined name
# (and any other sheet-qualified references) pointing at the renamed sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "template"

# Renaming the sheet does not retarget the workbook-level Print_Area
# defined name automatically, so restate it explicitly against the new
# sheet name (same absolute range as before: A1:M25).
$ws.PageSetup.PrintArea = "`$A`$1:`$M`$25"
